# Generate Report for Handoff
# Update the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps for the
# 7ad558eb-8d63-4656-a9ea-32da7a63fd8e file across the Overview, zh-cn and de-de
# sheets to reflect a newly generated handoff report.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 6 corresponds to 7ad558eb-8d63-4656-a9ea-32da7a63fd8e.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-27 10:42:01"

# zh-cn sheet: row 6 corresponds to 7ad558eb-8d63-4656-a9ea-32da7a63fd8e.md
# Column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-27 10:41:56"

# de-de sheet: row 6 corresponds to 7ad558eb-8d63-4656-a9ea-32da7a63fd8e.md
# Column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-27 10:42:01"
